# Update leve-profit calculation columns (H-N) across multiple sheets
# per latest market data refresh from scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1202726.5
$ws.Range("J17").Value = 1202726.5
$ws.Range("L17").Value = 3608179.5
$ws.Range("N17").Value = -3608515.5

$ws.Range("H116").Value = 5749.5
$ws.Range("I116").Value = 4586.125
$ws.Range("K116").Value = 4586.125
$ws.Range("M116").Value = -1144.125

$ws.Range("H132").Value = 12457.281
$ws.Range("I132").Value = 13011.667
$ws.Range("J132").Value = 11398.909
$ws.Range("K132").Value = 39035.001
$ws.Range("L132").Value = 34196.727
$ws.Range("M132").Value = -36505.001
$ws.Range("N132").Value = -39256.727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 993996.4
$ws.Range("I61").Value = 3432.8572
$ws.Range("K61").Value = 3432.8572
$ws.Range("M61").Value = -3220.8572

$ws.Range("H74").Value = 24486.77
$ws.Range("I74").Value = 1869.6471
$ws.Range("K74").Value = 1869.6471
$ws.Range("M74").Value = -995.6470999999999

$ws.Range("H77").Value = 24486.77
$ws.Range("I77").Value = 1869.6471
$ws.Range("K77").Value = 9348.235499999999
$ws.Range("M77").Value = -4980.235499999999

$ws.Range("H136").Value = 993996.4
$ws.Range("I136").Value = 3432.8572
$ws.Range("K136").Value = 10298.5716
$ws.Range("M136").Value = -7748.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 37392.516
$ws.Range("I134").Value = 36568.285
$ws.Range("K134").Value = 109704.855
$ws.Range("M134").Value = -107169.855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 2150
$ws.Range("I48").Value = 2150
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 2150
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -1674
$ws.Range("N48").ClearContents()

$ws.Range("H58").Value = 12421.889
$ws.Range("I58").Value = 4670.769
$ws.Range("K58").Value = 4670.769
$ws.Range("M58").Value = -4467.769

$ws.Range("H94").Value = 1127.75
$ws.Range("I94").Value = 999
$ws.Range("J94").Value = 1170.6666
$ws.Range("K94").Value = 999
$ws.Range("L94").Value = 1170.6666
$ws.Range("M94").Value = -548
$ws.Range("N94").Value = -2072.6666

$ws.Range("H129").Value = 92166.5
$ws.Range("J129").Value = 92166.5
$ws.Range("L129").Value = 92166.5
$ws.Range("N129").Value = -102166.5

$ws.Range("H134").Value = 28578064
$ws.Range("J134").Value = 100013720
$ws.Range("L134").Value = 300041160
$ws.Range("N134").Value = -300046230

$ws.Range("H136").Value = 12421.889
$ws.Range("I136").Value = 4670.769
$ws.Range("K136").Value = 14012.307
$ws.Range("M136").Value = -11462.307

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 148.81818
$ws.Range("I2").Value = 173.35715
$ws.Range("K2").Value = 1040.1429
$ws.Range("M2").Value = -927.1428999999998

$ws.Range("H3").Value = 1164.25
$ws.Range("I3").Value = 1164.25
$ws.Range("K3").Value = 3492.75
$ws.Range("M3").Value = -3380.75

$ws.Range("H7").Value = 6095.8335
$ws.Range("I7").Value = 131.25
$ws.Range("J7").Value = 7800
$ws.Range("K7").Value = 393.75
$ws.Range("L7").Value = 23400
$ws.Range("M7").Value = -281.75
$ws.Range("N7").Value = -23624

$ws.Range("H12").Value = 74.818184
$ws.Range("J12").Value = 95.625
$ws.Range("L12").Value = 286.875
$ws.Range("N12").Value = -632.875

$ws.Range("H29").Value = 2237.5454
$ws.Range("I29").Value = 2904.75
$ws.Range("J29").Value = 1856.2858
$ws.Range("K29").Value = 8714.25
$ws.Range("L29").Value = 5568.857400000001
$ws.Range("M29").Value = -8437.25
$ws.Range("N29").Value = -6122.857400000001

$ws.Range("H31").Value = 4200.25
$ws.Range("I31").Value = 4933
$ws.Range("J31").Value = 2002
$ws.Range("K31").Value = 14799
$ws.Range("L31").Value = 6006
$ws.Range("M31").Value = -14511
$ws.Range("N31").Value = -6582

$ws.Range("H38").Value = 157.0625
$ws.Range("I38").Value = 174.71428
$ws.Range("J38").Value = 143.33333
$ws.Range("K38").Value = 524.14284
$ws.Range("L38").Value = 429.99999
$ws.Range("M38").Value = -177.14284
$ws.Range("N38").Value = -1123.99999

$ws.Range("H112").Value = 8487.25
$ws.Range("J112").Value = 9650
$ws.Range("L112").Value = 28950
$ws.Range("N112").Value = -31166

$ws.Range("H113").Value = 1175.9412
$ws.Range("J113").Value = 1162.8182
$ws.Range("L113").Value = 3488.4546
$ws.Range("N113").Value = -7828.4546

$ws.Range("H133").Value = 7485.7144
$ws.Range("I133").Value = 4500
$ws.Range("J133").Value = 8680
$ws.Range("K133").Value = 13500
$ws.Range("L133").Value = 26040
$ws.Range("M133").Value = -8440
$ws.Range("N133").Value = -36160

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H126").Value = 7690.4546
$ws.Range("I126").Value = 10173.143
$ws.Range("J126").Value = 3345.75
$ws.Range("K126").Value = 30519.429
$ws.Range("L126").Value = 10037.25
$ws.Range("M126").Value = -28049.429
$ws.Range("N126").Value = -14977.25

$ws.Range("H129").Value = 35000
$ws.Range("J129").Value = 35000
$ws.Range("L129").Value = 35000
$ws.Range("N129").Value = -45000

$ws.Range("H132").Value = 778217.5
$ws.Range("I132").Value = 3771.2942
$ws.Range("K132").Value = 11313.8826
$ws.Range("M132").Value = -8783.882599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4333.3335
$ws.Range("I40").Value = 4000
$ws.Range("K40").Value = 4000
$ws.Range("M40").Value = -3864

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 18333.166
$ws.Range("I51").Value = 4999.5
$ws.Range("J51").Value = 25000
$ws.Range("K51").Value = 4999.5
$ws.Range("L51").Value = 25000
$ws.Range("M51").Value = -4489.5
$ws.Range("N51").Value = -26020

$ws.Range("H107").Value = 1322
$ws.Range("I107").Value = 1629.2222
$ws.Range("J107").Value = 400.33334
$ws.Range("K107").Value = 4887.6666
$ws.Range("L107").Value = 1201.00002
$ws.Range("M107").Value = -2967.6666
$ws.Range("N107").Value = -5041.000019999999

$ws.Range("H126").Value = 4261.9653
$ws.Range("I126").Value = 4601.1904
$ws.Range("K126").Value = 13803.5712
$ws.Range("M126").Value = -11333.5712

$ws.Range("H136").Value = 329698.44
$ws.Range("I136").Value = 1685.7391
$ws.Range("K136").Value = 5057.2173
$ws.Range("M136").Value = -2507.2173
